# Fix sorting and generate viable xlsx and charts
#
# The "Data" sheet holds benchmark results for RandomInt33 - HeapSortTimes.csv.
# The Avg_Time_ms column (D) for the two smallest row counts (5000 and 10000
# rows, in rows 2 and 3) was recomputed after fixing the sort, so update
# those two cells to the corrected averages. (The embedded scatter chart
# reads its series from Data!$D$2:$D$8 / Data!$E$2:$E$8, so it will reflect
# these corrected values too.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("D2").Value = 0.9105976200000001
$ws.Range("D3").Value = 2.0052926
